$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values as text (inline strings) in the
# original workbook. Force each updated Price cell to Text format first so that
# assigning a numeric-looking string (e.g. "325.78") does not get auto-converted
# into a real number by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.288.26"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.00"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.78"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4367"
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.76"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07668"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.141"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.302"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.493"
$ws.Range("E15").Value = "  +3.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.820.92"
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.95"
$ws.Range("E17").Value = "  +7.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001078"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06484"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.36"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.236"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.295.14"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.123"
$ws.Range("E25").Value = "  -8.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.38"
$ws.Range("E26").Value = "  +5.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.71"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.028.65"
$ws.Range("E28").Value = "  +3.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.270"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.19"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.209"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("E32").Value = "  +5.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09133"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.557"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.97"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.232"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06212"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.194"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.063"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.427"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.81"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6099"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.736"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.43"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.013"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.158"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06993"
$ws.Range("E51").Value = "  +1.34%  "
